$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 762.25
$ws.Range("I38").Value = 168.16667
$ws.Range("J38").Value = 2544.5
$ws.Range("K38").Value = 504.50001
$ws.Range("L38").Value = 7633.5
$ws.Range("M38").Value = -132.50001
$ws.Range("N38").Value = -8377.5

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2599.5715
$ws.Range("I45").Value = 1449.25
$ws.Range("J45").Value = 4133.3335
$ws.Range("K45").Value = 1449.25
$ws.Range("L45").Value = 4133.3335
$ws.Range("M45").Value = -1072.25
$ws.Range("N45").Value = -4887.3335

# Row 110
$ws.Range("H110").Value = 2346.12
$ws.Range("I110").Value = 1375.2142
$ws.Range("J110").Value = 3581.818
$ws.Range("K110").Value = 1375.2142
$ws.Range("L110").Value = 3581.818
$ws.Range("M110").Value = 669.7858000000001
$ws.Range("N110").Value = -7671.818

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1864.6666
$ws.Range("I86").Value = 1522.25
$ws.Range("J86").Value = 2549.5
$ws.Range("K86").Value = 1522.25
$ws.Range("L86").Value = 2549.5
$ws.Range("M86").Value = -399.25
$ws.Range("N86").Value = -4795.5

# Row 89
$ws.Range("H89").Value = 1864.6666
$ws.Range("I89").Value = 1522.25
$ws.Range("J89").Value = 2549.5
$ws.Range("K89").Value = 7611.25
$ws.Range("L89").Value = 12747.5
$ws.Range("M89").Value = -1995.25
$ws.Range("N89").Value = -23979.5

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 6469.6875
$ws.Range("I3").Value = 4497.1113
$ws.Range("J3").Value = 9005.857
$ws.Range("K3").Value = 4497.1113
$ws.Range("L3").Value = 9005.857
$ws.Range("M3").Value = -4384.1113
$ws.Range("N3").Value = -9231.857

# Row 16
$ws.Range("H16").Value = 2963.6667
$ws.Range("I16").Value = 2576.5
$ws.Range("J16").Value = 4899.5
$ws.Range("K16").Value = 2576.5
$ws.Range("L16").Value = 4899.5
$ws.Range("M16").Value = -2289.5
$ws.Range("N16").Value = -5473.5

# Row 22
$ws.Range("H22").Value = 84737.92
$ws.Range("I22").Value = 83466
$ws.Range("J22").Value = 100001
$ws.Range("K22").Value = 83466
$ws.Range("L22").Value = 100001
$ws.Range("M22").Value = -83116
$ws.Range("N22").Value = -100701

# Row 58
$ws.Range("H58").Value = 4997.5
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 4997.5
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 4997.5
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -5403.5

# Row 93
$ws.Range("H93").Value = 10253
$ws.Range("I93").Value = 10253
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 10253
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -8381

# Row 99
$ws.Range("H99").Value = 3102.1538
$ws.Range("I99").Value = 3171.9
$ws.Range("J99").Value = 2869.6667
$ws.Range("K99").Value = 3171.9
$ws.Range("L99").Value = 2869.6667
$ws.Range("M99").Value = -1673.9
$ws.Range("N99").Value = -5865.6667

# Row 113
$ws.Range("H113").Value = 2963.6667
$ws.Range("I113").Value = 2576.5
$ws.Range("J113").Value = 4899.5
$ws.Range("K113").Value = 2576.5
$ws.Range("L113").Value = 4899.5
$ws.Range("M113").Value = -406.5
$ws.Range("N113").Value = -9239.5

# Row 126
$ws.Range("H126").Value = 3102.1538
$ws.Range("I126").Value = 3171.9
$ws.Range("J126").Value = 2869.6667
$ws.Range("K126").Value = 9515.700000000001
$ws.Range("L126").Value = 8609.000100000001
$ws.Range("M126").Value = -7045.700000000001
$ws.Range("N126").Value = -13549.0001

# Row 136
$ws.Range("H136").Value = 4997.5
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 4997.5
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 14992.5
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -20092.5

$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()

# Row 113
$ws.Range("H113").Value = 866.6667
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 2400
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -230
$ws.Range("N113").Value = -7340

# Row 139
$ws.Range("H139").Value = 1803.4546
$ws.Range("I139").Value = 1803.4546
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 5410.3638
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -270.3638000000001
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 7
$ws.Range("H7").Value = 206059.8
$ws.Range("I7").Value = 1000000
$ws.Range("J7").Value = 7574.75
$ws.Range("K7").Value = 1000000
$ws.Range("L7").Value = 7574.75
$ws.Range("M7").Value = -999888
$ws.Range("N7").Value = -7798.75

# Row 8
$ws.Range("H8").Value = 206059.8
$ws.Range("I8").Value = 1000000
$ws.Range("J8").Value = 7574.75
$ws.Range("K8").Value = 1000000
$ws.Range("L8").Value = 7574.75
$ws.Range("M8").Value = -999861
$ws.Range("N8").Value = -7852.75

# Row 58
$ws.Range("H58").Value = 23166.334
$ws.Range("I58").Value = 23166.334
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 23166.334
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -22889.334
$ws.Range("N58").ClearContents()

# Row 131
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()

# Row 4
$ws.Range("H4").Value = 3858
$ws.Range("I4").Value = 3709
$ws.Range("J4").Value = 4007
$ws.Range("K4").Value = 3709
$ws.Range("L4").Value = 4007
$ws.Range("M4").Value = -3596
$ws.Range("N4").Value = -4233

# Row 14
$ws.Range("H14").Value = 19001.666
$ws.Range("I14").Value = 7500
$ws.Range("J14").Value = 24752.5
$ws.Range("K14").Value = 7500
$ws.Range("L14").Value = 24752.5
$ws.Range("M14").Value = -7328
$ws.Range("N14").Value = -25096.5

# Row 15
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()

# Row 18
$ws.Range("H18").Value = 28250
$ws.Range("I18").Value = 48000
$ws.Range("J18").Value = 8500
$ws.Range("K18").Value = 48000
$ws.Range("L18").Value = 8500
$ws.Range("M18").Value = -47828
$ws.Range("N18").Value = -8844

# Row 20
$ws.Range("H20").Value = 4666.6665
$ws.Range("I20").Value = 4000
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 4000
$ws.Range("L20").Value = 5000
$ws.Range("M20").Value = -3774
$ws.Range("N20").Value = -5452

# Row 28
$ws.Range("H28").Value = 3858
$ws.Range("I28").Value = 3709
$ws.Range("J28").Value = 4007
$ws.Range("K28").Value = 3709
$ws.Range("L28").Value = 4007
$ws.Range("M28").Value = -3477
$ws.Range("N28").Value = -4471

# Row 37
$ws.Range("H37").Value = 3858
$ws.Range("I37").Value = 3709
$ws.Range("J37").Value = 4007
$ws.Range("K37").Value = 3709
$ws.Range("L37").Value = 4007
$ws.Range("M37").Value = -3602
$ws.Range("N37").Value = -4221

# Row 132
$ws.Range("H132").Value = 3627.7144
$ws.Range("I132").Value = 3279.4
$ws.Range("J132").Value = 4498.5
$ws.Range("K132").Value = 9838.200000000001
$ws.Range("L132").Value = 13495.5
$ws.Range("M132").Value = -7308.200000000001
$ws.Range("N132").Value = -18555.5

$ws = $wb.Worksheets.Item("WVR")
# Row 11
$ws.Range("H11").Value = 1848164
$ws.Range("I11").Value = 5508500
$ws.Range("J11").Value = 17996
$ws.Range("K11").Value = 5508500
$ws.Range("L11").Value = 17996
$ws.Range("M11").Value = -5508358
$ws.Range("N11").Value = -18280

# Row 12
$ws.Range("H12").Value = 2899
$ws.Range("I12").Value = 2899
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 2899
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -2757
$ws.Range("N12").ClearContents()

# Row 20
$ws.Range("H20").Value = 15298.1
$ws.Range("I20").Value = 39995
$ws.Range("J20").Value = 12554
$ws.Range("K20").Value = 39995
$ws.Range("L20").Value = 12554
$ws.Range("M20").Value = -39755
$ws.Range("N20").Value = -13034

# Row 47
$ws.Range("H47").Value = 44997.5
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 44997.5
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 44997.5
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -46141.5

# Row 126
$ws.Range("H126").Value = 1949.5
$ws.Range("I126").Value = 1949.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5848.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3378.5
